$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update cryptocurrency price (column D) and volume-change (column E) values.
# NumberFormat is forced to text ("@") before assignment so that numeric-looking
# strings (e.g. "1.005") are stored as text, matching the source data which is
# always text (inline strings). ClearFormats() afterwards removes the temporary
# number-format override so the cell keeps its original (unstyled) appearance.

$cell = $ws.Range('D2')
$cell.NumberFormat = '@'
$cell.Value = '27.230.27'
$cell.ClearFormats()

$cell = $ws.Range('E2')
$cell.NumberFormat = '@'
$cell.Value = '  -1.79%  '
$cell.ClearFormats()

$cell = $ws.Range('D3')
$cell.NumberFormat = '@'
$cell.Value = '1.819.88'
$cell.ClearFormats()

$cell = $ws.Range('E3')
$cell.NumberFormat = '@'
$cell.Value = '  -2.19%  '
$cell.ClearFormats()

$cell = $ws.Range('D4')
$cell.NumberFormat = '@'
$cell.Value = '1.005'
$cell.ClearFormats()

$cell = $ws.Range('E4')
$cell.NumberFormat = '@'
$cell.Value = '  -1.50%  '
$cell.ClearFormats()

$cell = $ws.Range('D5')
$cell.NumberFormat = '@'
$cell.Value = '314.32'
$cell.ClearFormats()

$cell = $ws.Range('E5')
$cell.NumberFormat = '@'
$cell.Value = '  -2.12%  '
$cell.ClearFormats()

$cell = $ws.Range('D6')
$cell.NumberFormat = '@'
$cell.Value = '1.004'
$cell.ClearFormats()

$cell = $ws.Range('E6')
$cell.NumberFormat = '@'
$cell.Value = '  -1.50%  '
$cell.ClearFormats()

$cell = $ws.Range('D7')
$cell.NumberFormat = '@'
$cell.Value = '0.4267'
$cell.ClearFormats()

$cell = $ws.Range('E7')
$cell.NumberFormat = '@'
$cell.Value = '  -2.45%  '
$cell.ClearFormats()

$cell = $ws.Range('D8')
$cell.NumberFormat = '@'
$cell.Value = '0.3669'
$cell.ClearFormats()

$cell = $ws.Range('E8')
$cell.NumberFormat = '@'
$cell.Value = '  -2.95%  '
$cell.ClearFormats()

$cell = $ws.Range('D9')
$cell.NumberFormat = '@'
$cell.Value = '0.07215'
$cell.ClearFormats()

$cell = $ws.Range('E9')
$cell.NumberFormat = '@'
$cell.Value = '  -2.69%  '
$cell.ClearFormats()

$cell = $ws.Range('D10')
$cell.NumberFormat = '@'
$cell.Value = '0.8599'
$cell.ClearFormats()

$cell = $ws.Range('E10')
$cell.NumberFormat = '@'
$cell.Value = '  -2.81%  '
$cell.ClearFormats()

$cell = $ws.Range('E11')
$cell.NumberFormat = '@'
$cell.Value = '  -3.16%  '
$cell.ClearFormats()

$cell = $ws.Range('D12')
$cell.NumberFormat = '@'
$cell.Value = '1.828.39'
$cell.ClearFormats()

$cell = $ws.Range('E12')
$cell.NumberFormat = '@'
$cell.Value = '  -1.83%  '
$cell.ClearFormats()

$cell = $ws.Range('D13')
$cell.NumberFormat = '@'
$cell.Value = '6.663'
$cell.ClearFormats()

$cell = $ws.Range('E13')
$cell.NumberFormat = '@'
$cell.Value = '  -1.47%  '
$cell.ClearFormats()

$cell = $ws.Range('D14')
$cell.NumberFormat = '@'
$cell.Value = '0.07108'
$cell.ClearFormats()

$cell = $ws.Range('E14')
$cell.NumberFormat = '@'
$cell.Value = '  -0.58%  '
$cell.ClearFormats()

$cell = $ws.Range('D15')
$cell.NumberFormat = '@'
$cell.Value = '5.298'
$cell.ClearFormats()

$cell = $ws.Range('E15')
$cell.NumberFormat = '@'
$cell.Value = '  -3.71%  '
$cell.ClearFormats()

$cell = $ws.Range('D16')
$cell.NumberFormat = '@'
$cell.Value = '87.96'
$cell.ClearFormats()

$cell = $ws.Range('E16')
$cell.NumberFormat = '@'
$cell.Value = '  +0.19%  '
$cell.ClearFormats()

$cell = $ws.Range('E17')
$cell.NumberFormat = '@'
$cell.Value = '  -1.75%  '
$cell.ClearFormats()

$cell = $ws.Range('D18')
$cell.NumberFormat = '@'
$cell.Value = '0.000008861'
$cell.ClearFormats()

$cell = $ws.Range('E18')
$cell.NumberFormat = '@'
$cell.Value = '  -2.08%  '
$cell.ClearFormats()

$cell = $ws.Range('E19')
$cell.NumberFormat = '@'
$cell.Value = '  -1.52%  '
$cell.ClearFormats()

$cell = $ws.Range('E20')
$cell.NumberFormat = '@'
$cell.Value = '  -2.95%  '
$cell.ClearFormats()

$cell = $ws.Range('D21')
$cell.NumberFormat = '@'
$cell.Value = '27.246.28'
$cell.ClearFormats()

$cell = $ws.Range('E21')
$cell.NumberFormat = '@'
$cell.Value = '  -1.74%  '
$cell.ClearFormats()

$cell = $ws.Range('D22')
$cell.NumberFormat = '@'
$cell.Value = '5.142'
$cell.ClearFormats()

$cell = $ws.Range('E23')
$cell.NumberFormat = '@'
$cell.Value = '  -2.70%  '
$cell.ClearFormats()

$cell = $ws.Range('D24')
$cell.NumberFormat = '@'
$cell.Value = '2.044.48'
$cell.ClearFormats()

$cell = $ws.Range('E24')
$cell.NumberFormat = '@'
$cell.Value = '  -2.37%  '
$cell.ClearFormats()

$cell = $ws.Range('D25')
$cell.NumberFormat = '@'
$cell.Value = '2.004'
$cell.ClearFormats()

$cell = $ws.Range('E25')
$cell.NumberFormat = '@'
$cell.Value = '  -1.75%  '
$cell.ClearFormats()

$cell = $ws.Range('D26')
$cell.NumberFormat = '@'
$cell.Value = '153.10'
$cell.ClearFormats()

$cell = $ws.Range('E26')
$cell.NumberFormat = '@'
$cell.Value = '  -2.68%  '
$cell.ClearFormats()

$cell = $ws.Range('D27')
$cell.NumberFormat = '@'
$cell.Value = '18.28'
$cell.ClearFormats()

$cell = $ws.Range('E27')
$cell.NumberFormat = '@'
$cell.Value = '  -2.33%  '
$cell.ClearFormats()

$cell = $ws.Range('D28')
$cell.NumberFormat = '@'
$cell.Value = '2.115'
$cell.ClearFormats()

$cell = $ws.Range('E28')
$cell.NumberFormat = '@'
$cell.Value = '  +6.26%  '
$cell.ClearFormats()

$cell = $ws.Range('E29')
$cell.NumberFormat = '@'
$cell.Value = '  -3.90%  '
$cell.ClearFormats()

$cell = $ws.Range('D30')
$cell.NumberFormat = '@'
$cell.Value = '116.13'
$cell.ClearFormats()

$cell = $ws.Range('E30')
$cell.NumberFormat = '@'
$cell.Value = '  -4.48%  '
$cell.ClearFormats()

$cell = $ws.Range('D31')
$cell.NumberFormat = '@'
$cell.Value = '0.08887'
$cell.ClearFormats()

$cell = $ws.Range('E31')
$cell.NumberFormat = '@'
$cell.Value = '  -1.93%  '
$cell.ClearFormats()

$cell = $ws.Range('D32')
$cell.NumberFormat = '@'
$cell.Value = '1.193'
$cell.ClearFormats()

$cell = $ws.Range('E32')
$cell.NumberFormat = '@'
$cell.Value = '  -1.93%  '
$cell.ClearFormats()

$cell = $ws.Range('D33')
$cell.NumberFormat = '@'
$cell.Value = '0.7581'
$cell.ClearFormats()

$cell = $ws.Range('E33')
$cell.NumberFormat = '@'
$cell.Value = '  -1.50%  '
$cell.ClearFormats()

$cell = $ws.Range('D34')
$cell.NumberFormat = '@'
$cell.Value = '4.449'
$cell.ClearFormats()

$cell = $ws.Range('D35')
$cell.NumberFormat = '@'
$cell.Value = '2.819'
$cell.ClearFormats()

$cell = $ws.Range('E35')
$cell.NumberFormat = '@'
$cell.Value = '  -7.14%  '
$cell.ClearFormats()

$cell = $ws.Range('D36')
$cell.NumberFormat = '@'
$cell.Value = '1.004'
$cell.ClearFormats()

$cell = $ws.Range('E36')
$cell.NumberFormat = '@'
$cell.Value = '  -1.61%  '
$cell.ClearFormats()

$cell = $ws.Range('D37')
$cell.NumberFormat = '@'
$cell.Value = '1.113'
$cell.ClearFormats()

$cell = $ws.Range('E37')
$cell.NumberFormat = '@'
$cell.Value = '  -2.32%  '
$cell.ClearFormats()

$cell = $ws.Range('D38')
$cell.NumberFormat = '@'
$cell.Value = '0.01961'
$cell.ClearFormats()

$cell = $ws.Range('E38')
$cell.NumberFormat = '@'
$cell.Value = '  -0.91%  '
$cell.ClearFormats()

$cell = $ws.Range('D39')
$cell.NumberFormat = '@'
$cell.Value = '0.05261'
$cell.ClearFormats()

$cell = $ws.Range('D40')
$cell.NumberFormat = '@'
$cell.Value = '2.900'
$cell.ClearFormats()

$cell = $ws.Range('E40')
$cell.NumberFormat = '@'
$cell.Value = '  +1.05%  '
$cell.ClearFormats()

$cell = $ws.Range('D41')
$cell.NumberFormat = '@'
$cell.Value = '7.121'
$cell.ClearFormats()

$cell = $ws.Range('E41')
$cell.NumberFormat = '@'
$cell.Value = '  +2.19%  '
$cell.ClearFormats()

$cell = $ws.Range('D42')
$cell.NumberFormat = '@'
$cell.Value = '0.1676'
$cell.ClearFormats()

$cell = $ws.Range('E42')
$cell.NumberFormat = '@'
$cell.Value = '  -0.23%  '
$cell.ClearFormats()

$cell = $ws.Range('D43')
$cell.NumberFormat = '@'
$cell.Value = '0.5018'
$cell.ClearFormats()

$cell = $ws.Range('D44')
$cell.NumberFormat = '@'
$cell.Value = '8.607'
$cell.ClearFormats()

$cell = $ws.Range('E44')
$cell.NumberFormat = '@'
$cell.Value = '  -1.26%  '
$cell.ClearFormats()

$cell = $ws.Range('D45')
$cell.NumberFormat = '@'
$cell.Value = '10.64'
$cell.ClearFormats()

$cell = $ws.Range('E45')
$cell.NumberFormat = '@'
$cell.Value = '  -1.51%  '
$cell.ClearFormats()

$cell = $ws.Range('D46')
$cell.NumberFormat = '@'
$cell.Value = '106.57'
$cell.ClearFormats()

$cell = $ws.Range('E46')
$cell.NumberFormat = '@'
$cell.Value = '  -3.40%  '
$cell.ClearFormats()

$cell = $ws.Range('D47')
$cell.NumberFormat = '@'
$cell.Value = '0.4699'
$cell.ClearFormats()

$cell = $ws.Range('E47')
$cell.NumberFormat = '@'
$cell.Value = '  -0.69%  '
$cell.ClearFormats()

$cell = $ws.Range('D48')
$cell.NumberFormat = '@'
$cell.Value = '1.004'
$cell.ClearFormats()

$cell = $ws.Range('E48')
$cell.NumberFormat = '@'
$cell.Value = '  -1.63%  '
$cell.ClearFormats()

$cell = $ws.Range('D49')
$cell.NumberFormat = '@'
$cell.Value = '0.06385'
$cell.ClearFormats()

$cell = $ws.Range('E49')
$cell.NumberFormat = '@'
$cell.Value = '  -1.54%  '
$cell.ClearFormats()

$cell = $ws.Range('D50')
$cell.NumberFormat = '@'
$cell.Value = '1.660'
$cell.ClearFormats()

$cell = $ws.Range('E50')
$cell.NumberFormat = '@'
$cell.Value = '  -3.15%  '
$cell.ClearFormats()

$cell = $ws.Range('D51')
$cell.NumberFormat = '@'
$cell.Value = '1.795'
$cell.ClearFormats()

$cell = $ws.Range('E51')
$cell.NumberFormat = '@'
$cell.Value = '  -2.87%  '
$cell.ClearFormats()

